$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the duplicate/erroneous "shy US EQUITY" row (was row 9).
#    Everything below shifts up by one.
# ------------------------------------------------------------------
$ws.Rows(9).Delete()

# ------------------------------------------------------------------
# 2. Add the new "NomeCurto" (short name) column C, writing the cells
#    in the same order the source workbook was authored in so the
#    shared-string table comes out in the matching sequence.
# ------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "NomeCurto"
$ws.Cells.Item(3, 3).Value = "Fortune FI Fund"
$ws.Cells.Item(2, 3).Value = "US T-Bills"
$ws.Cells.Item(4, 3).Value = "CoCos AT1"
$ws.Cells.Item(5, 3).Value = "US TIPS"
$ws.Cells.Item(6, 3).Value = "US Treas 7–10y"
$ws.Cells.Item(7, 3).Value = "US Treas 20y+"
$ws.Cells.Item(8, 3).Value = "EM Sov Bonds"
$ws.Cells.Item(12, 3).Value = "US Treas 1–3y"
$ws.Cells.Item(9, 3).Value = "US HY Credit"
$ws.Cells.Item(10, 3).Value = "EM Corp Bonds"
$ws.Cells.Item(11, 3).Value = "US IG Credit"
$ws.Cells.Item(13, 3).Value = "Nasdaq 100"
$ws.Cells.Item(14, 3).Value = "S&P 500"
$ws.Cells.Item(15, 3).Value = "Russell 2000"
$ws.Cells.Item(16, 3).Value = "S&P 500 EqW"
$ws.Cells.Item(17, 3).Value = "Cloud ETF"
$ws.Cells.Item(18, 3).Value = "Semis ETF"
$ws.Cells.Item(19, 3).Value = "Defense ETF"
$ws.Cells.Item(20, 3).Value = "Gold Miners"
$ws.Cells.Item(21, 3).Value = "Metals & Mining ETF"
$ws.Cells.Item(22, 3).Value = "Brazil Eq ETF"
$ws.Cells.Item(23, 3).Value = "China Large Cap ETF"
$ws.Cells.Item(24, 3).Value = "Mexico Eq ETF"
$ws.Cells.Item(25, 3).Value = "China Internet ETF"
$ws.Cells.Item(26, 3).Value = "Japan Eq ETF"
$ws.Cells.Item(27, 3).Value = "Korea Eq ETF"
$ws.Cells.Item(28, 3).Value = "Europe Eq ETF"
$ws.Cells.Item(29, 3).Value = "Oil ETF"
$ws.Cells.Item(30, 3).Value = "Gold ETF"

# ------------------------------------------------------------------
# 3. Formatting.
#    Use a scratch cell + Copy/PasteSpecial(Formats) so each target
#    range picks up exactly one new cell style (matching the author's
#    two extra cellXfs entries) instead of accumulating an
#    intermediate style per property set.
# ------------------------------------------------------------------

# 3a. Column C data rows (2-30): wrap text + vertically centered.
$wrapScratch = $ws.Range("Z1")
$wrapScratch.WrapText = $true
$wrapScratch.VerticalAlignment = -4108
$wrapScratch.Copy()
$ws.Range("C2:C30").PasteSpecial(-4122)
$wrapScratch.Clear()

# 3b. Header row (A1:C1): bold.
$boldScratch = $ws.Range("Z2")
$boldScratch.Font.Bold = $true
$boldScratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$boldScratch.Clear()

# ------------------------------------------------------------------
# 4. Tidy up the view: drop the stray active-cell selection left over
#    from the previous session.
# ------------------------------------------------------------------
$ws.Range("A1").Select()
